$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.743.89"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3
$ws.Range("D3").Value = "3.468.39"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'577.21"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6
$ws.Range("D6").Value = "'145.29"
$ws.Range("E6").Value = "  -2.47%  "

# Row 7
$ws.Range("D7").Value = "3.474.44"
$ws.Range("E7").Value = "  +0.52%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.481"
$ws.Range("E9").Value = "  +1.20%  "

# Row 10
$ws.Range("D10").Value = "'7.64"
$ws.Range("E10").Value = "  -0.73%  "

# Row 11
$ws.Range("E11").Value = "  +2.51%  "

# Row 12
$ws.Range("D12").Value = "'0.392"
$ws.Range("E12").Value = "  +0.05%  "

# Row 13
$ws.Range("D13").Value = "4.041.05"
$ws.Range("E13").Value = "  -0.12%  "

# Row 14
$ws.Range("D14").Value = "'28.93"
$ws.Range("E14").Value = "  +5.71%  "

# Row 15
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("E16").Value = "  -0.41%  "

# Row 17
$ws.Range("D17").Value = "3.452.08"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("D18").Value = "61.776.20"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "  +3.81%  "

# Row 20
$ws.Range("D20").Value = "'14.49"
$ws.Range("E20").Value = "  +2.51%  "

# Row 21
$ws.Range("D21").Value = "'9.49"
$ws.Range("E21").Value = "  -0.84%  "

# Row 22
$ws.Range("D22").Value = "'400.15"
$ws.Range("E22").Value = "  +4.48%  "

# Row 23
$ws.Range("D23").Value = "'0.568"
$ws.Range("E23").Value = "  +0.93%  "

# Row 24
$ws.Range("D24").Value = "'74.18"
$ws.Range("E24").Value = "  +2.63%  "

# Row 25
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.28%  "

# Row 26
$ws.Range("D26").Value = "'0.0000124"
$ws.Range("E26").Value = "  -0.69%  "

# Row 27
$ws.Range("D27").Value = "3.584.22"
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("D28").Value = "'0.180"
$ws.Range("E28").Value = "  +1.44%  "

# Row 29
$ws.Range("D29").Value = "'7.71"
$ws.Range("E29").Value = "  -1.43%  "

# Row 30
$ws.Range("E30").Value = "  -0.10%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.29"
$ws.Range("E31").Value = "  +0.29%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.49"
$ws.Range("E32").Value = "  -6.91%  "

# Row 33
$ws.Range("D33").Value = "'2.19"
$ws.Range("E33").Value = "  +0.87%  "

# Row 34
$ws.Range("E34").Value = "  -0.08%  "

# Row 35
$ws.Range("D35").Value = "'24.11"
$ws.Range("E35").Value = "  +0.40%  "

# Row 36
$ws.Range("D36").Value = "'7.10"
$ws.Range("E36").Value = "  +0.97%  "

# Row 37
$ws.Range("D37").Value = "3.486.63"
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("D38").Value = "'5.21"
$ws.Range("E38").Value = "  -2.03%  "

# Row 39
$ws.Range("D39").Value = "'1.57"
$ws.Range("E39").Value = "  -0.52%  "

# Row 40
$ws.Range("D40").Value = "'167.92"
$ws.Range("E40").Value = "  +0.56%  "

# Row 41
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'28.66"
$ws.Range("E41").Value = "  +9.15%  "

# Row 42
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0796"
$ws.Range("E42").Value = "  +0.96%  "

# Row 43
$ws.Range("D43").Value = "'0.807"
$ws.Range("E43").Value = "  +1.45%  "

# Row 44
$ws.Range("D44").Value = "'4.57"
$ws.Range("E44").Value = "  +2.38%  "

# Row 45
$ws.Range("D45").Value = "'1.75"
$ws.Range("E45").Value = "  +1.29%  "

# Row 46
$ws.Range("D46").Value = "'0.996"
$ws.Range("E46").Value = "  -0.52%  "

# Row 47
$ws.Range("D47").Value = "2.652.17"
$ws.Range("E47").Value = "  +0.05%  "

# Row 48
$ws.Range("E48").Value = "  -3.73%  "

# Row 49
$ws.Range("D49").Value = "'6.98"
$ws.Range("E49").Value = "  +1.26%  "

# Row 50
$ws.Range("D50").Value = "'23.08"
$ws.Range("E50").Value = "  -3.98%  "

# Row 51
$ws.Range("D51").Value = "'2.43"
$ws.Range("E51").Value = "  +3.91%  "
